# Updates the "df_route" test-journey sheet: rows 2-26 get new origin/destination
# address data (BA11 5LB -> BA11 5AP, Frome) replacing the old Midsomer
# Norton -> Henstridge data, and four additional rows (27-30) are appended
# with the same address/distance data and new Lat/Lng coordinates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origin = "BA11 5LB"
$destination = "BA11 5AP"
$startAddress = "81 Knights Maltings, Frome, Frome, BA11 5LB, United Kingdom"
$endAddress = "55 Tower View, Frome, Frome, BA11 5AP, United Kingdom"
$distanceText = 3.0501
$distanceValue = 3050.1
$durationText = 8.711666666666668
$durationValue = 522.7

# Lat/Lng per row (rows 2 through 30)
$latLng = @{
    2 = @(51.22234, -2.31109)
    3 = @(51.22237, -2.3107)
    4 = @(51.22273, -2.31064)
    5 = @(51.22283, -2.31005)
    6 = @(51.22298, -2.30982)
    7 = @(51.22374, -2.30909)
    8 = @(51.22498, -2.30754)
    9 = @(51.22534, -2.30686)
    10 = @(51.22581, -2.3054)
    11 = @(51.22681, -2.30373)
    12 = @(51.22708, -2.30363)
    13 = @(51.22726, -2.30377)
    14 = @(51.22884, -2.3063)
    15 = @(51.22893, -2.30699)
    16 = @(51.22876, -2.30829)
    17 = @(51.22791, -2.31099)
    18 = @(51.22768, -2.31325)
    19 = @(51.22642, -2.31437)
    20 = @(51.22582, -2.31544)
    21 = @(51.22519, -2.31769)
    22 = @(51.22484, -2.32133)
    23 = @(51.22421, -2.32144)
    24 = @(51.22301, -2.32124)
    25 = @(51.22183, -2.32125)
    26 = @(51.22048, -2.32074)
    27 = @(51.22034, -2.31956)
    28 = @(51.21988, -2.31827)
    29 = @(51.22045, -2.31728)
    30 = @(51.22032, -2.31717)
}

# First, rewrite the existing data rows (2-26) with the new address data.
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 2).Value = $origin
    $ws.Cells.Item($r, 3).Value = $destination
    $ws.Cells.Item($r, 4).Value = $startAddress
    $ws.Cells.Item($r, 5).Value = $endAddress
    $ws.Cells.Item($r, 6).Value = $distanceText
    $ws.Cells.Item($r, 7).Value = $distanceValue
    $ws.Cells.Item($r, 8).Value = $durationText
    $ws.Cells.Item($r, 9).Value = $durationValue
    $coords = $latLng[$r]
    $ws.Cells.Item($r, 10).Value = $coords[0]
    $ws.Cells.Item($r, 11).Value = $coords[1]
}

# Now append the four new rows (27-30), copying the formatting of the
# last existing data row (26) down so column A keeps its bold/centered
# bordered style.
$ws.Range("A26:K26").Copy()
$ws.Range("A27:K30").PasteSpecial(-4122)

for ($r = 27; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
    $ws.Cells.Item($r, 2).Value = $origin
    $ws.Cells.Item($r, 3).Value = $destination
    $ws.Cells.Item($r, 4).Value = $startAddress
    $ws.Cells.Item($r, 5).Value = $endAddress
    $ws.Cells.Item($r, 6).Value = $distanceText
    $ws.Cells.Item($r, 7).Value = $distanceValue
    $ws.Cells.Item($r, 8).Value = $durationText
    $ws.Cells.Item($r, 9).Value = $durationValue
    $coords = $latLng[$r]
    $ws.Cells.Item($r, 10).Value = $coords[0]
    $ws.Cells.Item($r, 11).Value = $coords[1]
}
